# Auto-generated edit script applying numeric corrections to the
# per-profession leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 9999
$ws.Range("J48").Value = 9999
$ws.Range("L48").Value = 29997
$ws.Range("N48").Value = -30581
$ws.Range("H51").Value = 113891950
$ws.Range("I51").Value = 205002830
$ws.Range("K51").Value = 205002830
$ws.Range("M51").Value = -205002346
$ws.Range("H56").Value = 9999
$ws.Range("J56").Value = 9999
$ws.Range("L56").Value = 29997
$ws.Range("N56").Value = -31065
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H86").Value = 5634.25
$ws.Range("I86").Value = 1600
$ws.Range("J86").Value = 6441.1
$ws.Range("K86").Value = 1600
$ws.Range("L86").Value = 6441.1
$ws.Range("M86").Value = -477
$ws.Range("N86").Value = -8687.1
$ws.Range("H89").Value = 5634.25
$ws.Range("I89").Value = 1600
$ws.Range("J89").Value = 6441.1
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 32205.5
$ws.Range("M89").Value = -2384
$ws.Range("N89").Value = -43437.5
$ws.Range("H103").Value = 498.66666
$ws.Range("I103").Value = 232.88889
$ws.Range("J103").Value = 698
$ws.Range("K103").Value = 698.6666700000001
$ws.Range("L103").Value = 2094
$ws.Range("M103").Value = -112.6666700000001
$ws.Range("N103").Value = -3266
$ws.Range("H127").Value = 3397.5
$ws.Range("I127").Value = 1295.5
$ws.Range("J127").Value = 5499.5
$ws.Range("K127").Value = 3886.5
$ws.Range("L127").Value = 16498.5
$ws.Range("M127").Value = 1073.5
$ws.Range("N127").Value = -26418.5
$ws.Range("H129").Value = 58826708
$ws.Range("I129").Value = 83333760
$ws.Range("J129").Value = 9783.200000000001
$ws.Range("K129").Value = 250001280
$ws.Range("L129").Value = 29349.6
$ws.Range("M129").Value = -249996280
$ws.Range("N129").Value = -39349.60000000001
$ws.Range("H135").Value = 1629.8462
$ws.Range("I135").Value = 513.0476
$ws.Range("K135").Value = 4617.4284
$ws.Range("M135").Value = -2082.4284
$ws.Range("H137").Value = 1838.0938
$ws.Range("I137").Value = 1909
$ws.Range("J137").Value = 1702.7273
$ws.Range("K137").Value = 5727
$ws.Range("L137").Value = 5108.1819
$ws.Range("M137").Value = -3177
$ws.Range("N137").Value = -10208.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8383.022999999999
$ws.Range("I32").Value = 8730.154
$ws.Range("K32").Value = 8730.154
$ws.Range("M32").Value = -8443.154
$ws.Range("H45").Value = 2561.5293
$ws.Range("J45").Value = 3339.8333
$ws.Range("L45").Value = 3339.8333
$ws.Range("N45").Value = -4093.8333
$ws.Range("H61").Value = 4607.797
$ws.Range("I61").Value = 3103.3408
$ws.Range("J61").Value = 7255.64
$ws.Range("K61").Value = 3103.3408
$ws.Range("L61").Value = 7255.64
$ws.Range("M61").Value = -2891.3408
$ws.Range("N61").Value = -7679.64
$ws.Range("H102").Value = 15388074
$ws.Range("I102").Value = 2449.25
$ws.Range("J102").Value = 40005076
$ws.Range("K102").Value = 2449.25
$ws.Range("L102").Value = 40005076
$ws.Range("M102").Value = -827.25
$ws.Range("N102").Value = -40008320
$ws.Range("H122").Value = 2534.5366
$ws.Range("I122").Value = 2229.818
$ws.Range("K122").Value = 6689.454000000001
$ws.Range("M122").Value = -4239.454000000001
$ws.Range("H132").Value = 2373.3833
$ws.Range("I132").Value = 2222.5557
$ws.Range("K132").Value = 6667.6671
$ws.Range("M132").Value = -4137.6671
$ws.Range("H136").Value = 4607.797
$ws.Range("I136").Value = 3103.3408
$ws.Range("J136").Value = 7255.64
$ws.Range("K136").Value = 9310.0224
$ws.Range("L136").Value = 21766.92
$ws.Range("M136").Value = -6760.0224
$ws.Range("N136").Value = -26866.92

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3711834.8
$ws.Range("I86").Value = 6072484
$ws.Range("K86").Value = 6072484
$ws.Range("M86").Value = -6071361
$ws.Range("H89").Value = 3711834.8
$ws.Range("I89").Value = 6072484
$ws.Range("K89").Value = 30362420
$ws.Range("M89").Value = -30356804
$ws.Range("H107").Value = 1871.7778
$ws.Range("I107").Value = 1786.25
$ws.Range("J107").Value = 1940.2
$ws.Range("K107").Value = 1786.25
$ws.Range("L107").Value = 1940.2
$ws.Range("M107").Value = 133.75
$ws.Range("N107").Value = -5780.2
$ws.Range("H134").Value = 9349.030000000001
$ws.Range("J134").Value = 9999.166999999999
$ws.Range("L134").Value = 29997.501
$ws.Range("N134").Value = -35067.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2008.9056
$ws.Range("I31").Value = 1207.5676
$ws.Range("K31").Value = 1207.5676
$ws.Range("M31").Value = -912.5676000000001
$ws.Range("H34").Value = 2008.9056
$ws.Range("I34").Value = 1207.5676
$ws.Range("K34").Value = 1207.5676
$ws.Range("M34").Value = -1005.5676
$ws.Range("H122").Value = 3230.577
$ws.Range("I122").Value = 3268.111
$ws.Range("K122").Value = 9804.332999999999
$ws.Range("M122").Value = -7354.332999999999
$ws.Range("H132").Value = 1483955
$ws.Range("J132").Value = 2495
$ws.Range("L132").Value = 7485
$ws.Range("N132").Value = -12545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 478.25
$ws.Range("I33").Value = 1018.5
$ws.Range("K33").Value = 6111
$ws.Range("M33").Value = -5828
$ws.Range("H47").Value = 3625
$ws.Range("I47").Value = 6500
$ws.Range("J47").Value = 750
$ws.Range("K47").Value = 19500
$ws.Range("L47").Value = 2250
$ws.Range("M47").Value = -19069
$ws.Range("N47").Value = -3112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 39015.715
$ws.Range("I99").Value = 30000
$ws.Range("J99").Value = 40518.332
$ws.Range("K99").Value = 30000
$ws.Range("L99").Value = 40518.332
$ws.Range("M99").Value = -27754
$ws.Range("N99").Value = -45010.332
$ws.Range("H126").Value = 2508.6
$ws.Range("I126").Value = 2137
$ws.Range("K126").Value = 6411
$ws.Range("M126").Value = -3941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4011
$ws.Range("I16").Value = 3444.5
$ws.Range("J16").Value = 4172.857
$ws.Range("K16").Value = 3444.5
$ws.Range("L16").Value = 4172.857
$ws.Range("M16").Value = -3274.5
$ws.Range("N16").Value = -4512.857
$ws.Range("H22").Value = 3429.9092
$ws.Range("I22").Value = 745
$ws.Range("J22").Value = 3698.4
$ws.Range("K22").Value = 745
$ws.Range("L22").Value = 3698.4
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -4288.4
$ws.Range("H27").Value = 3429.9092
$ws.Range("I27").Value = 745
$ws.Range("J27").Value = 3698.4
$ws.Range("K27").Value = 745
$ws.Range("L27").Value = 3698.4
$ws.Range("M27").Value = -638
$ws.Range("N27").Value = -3912.4
$ws.Range("H40").Value = 3810.6
$ws.Range("I40").Value = 3615
$ws.Range("J40").Value = 4267
$ws.Range("K40").Value = 3615
$ws.Range("L40").Value = 4267
$ws.Range("M40").Value = -3479
$ws.Range("N40").Value = -4539
$ws.Range("H68").Value = 3604.1667
$ws.Range("J68").Value = 3611
$ws.Range("L68").Value = 3611
$ws.Range("N68").Value = -5109
$ws.Range("H71").Value = 3604.1667
$ws.Range("J71").Value = 3611
$ws.Range("L71").Value = 18055
$ws.Range("N71").Value = -25543
$ws.Range("H93").Value = 2920.111
$ws.Range("I93").Value = 2276.2
$ws.Range("J93").Value = 3725
$ws.Range("K93").Value = 2276.2
$ws.Range("L93").Value = 3725
$ws.Range("M93").Value = -1028.2
$ws.Range("N93").Value = -6221
$ws.Range("H122").Value = 3017.8696
$ws.Range("I122").Value = 2320.65
$ws.Range("K122").Value = 6961.950000000001
$ws.Range("M122").Value = -4511.950000000001
$ws.Range("H132").Value = 3092.1277
$ws.Range("I132").Value = 3073.842
$ws.Range("K132").Value = 9221.526
$ws.Range("M132").Value = -6691.526
$ws.Range("H136").Value = 4744.1113
$ws.Range("I136").Value = 4471.3335
$ws.Range("J136").Value = 5698.8335
$ws.Range("K136").Value = 13414.0005
$ws.Range("L136").Value = 17096.5005
$ws.Range("M136").Value = -10864.0005
$ws.Range("N136").Value = -22196.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 198412.75
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 198412.75
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 48956.566
$ws.Range("I81").Value = 74189
$ws.Range("J81").Value = 9706.111000000001
$ws.Range("K81").Value = 148378
$ws.Range("L81").Value = 19412.222
$ws.Range("M81").Value = -147317
$ws.Range("N81").Value = -21534.222
$ws.Range("H84").Value = 48956.566
$ws.Range("I84").Value = 74189
$ws.Range("J84").Value = 9706.111000000001
$ws.Range("K84").Value = 741890
$ws.Range("L84").Value = 97061.11000000002
$ws.Range("M84").Value = -736586
$ws.Range("N84").Value = -107669.11
$ws.Range("H100").Value = 2930.25
$ws.Range("I100").Value = 2532.4285
$ws.Range("J100").Value = 3487.2
$ws.Range("K100").Value = 5064.857
$ws.Range("L100").Value = 6974.4
$ws.Range("M100").Value = -4523.857
$ws.Range("N100").Value = -8056.4
$ws.Range("H132").Value = 5661.778
$ws.Range("I132").Value = 8233
$ws.Range("J132").Value = 2447.75
$ws.Range("K132").Value = 24699
$ws.Range("L132").Value = 7343.25
$ws.Range("M132").Value = -22169
$ws.Range("N132").Value = -12403.25
